$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = 112154273
$ws.Range("Q9").Value = 442542
$ws.Range("R9").Value = 6204451
$ws.Range("AO9").Value = "på blmr av åkervädd m fl"
$ws.Range("A10").Value = 112145588
$ws.Range("B10").Value = 42584
$ws.Range("E10").Value = 100770
$ws.Range("F10").Value = "Mindre blåvinge"
$ws.Range("G10").Value = "Cupido minimus"
$ws.Range("H10").Value = "(Fuessly, 1775)"
$ws.Range("I10").Value = "'1"
$ws.Range("K10").Value = "imago/adult"
$ws.Range("L10").Value = "hane"
$ws.Range("M10").Value = "födosökande"
$ws.Range("Q10").Value = 442603
$ws.Range("R10").Value = 6204402
$ws.Range("Y10").Value = "2013-06-12"
$ws.Range("AA10").Value = "2013-06-12"
$ws.Range("AI10").Value = "i igenväxande hed"
$ws.Range("AO10").Value = "på praktveronika"
$ws.Range("AX10").Value = "Nils Otto Nilsson"
$ws.Range("A11").Value = 112154281
$ws.Range("B11").Value = 42584
$ws.Range("E11").Value = 100770
$ws.Range("F11").Value = "Mindre blåvinge"
$ws.Range("G11").Value = "Cupido minimus"
$ws.Range("H11").Value = "(Fuessly, 1775)"
$ws.Range("K11").ClearContents()
$ws.Range("L11").ClearContents()
$ws.Range("M11").Value = "friflygande"
$ws.Range("Q11").Value = 442664
$ws.Range("R11").Value = 6204260
$ws.Range("Y11").Value = "2013-07-09"
$ws.Range("AA11").Value = "2013-07-09"
$ws.Range("AC11").ClearContents()
$ws.Range("AI11").Value = "på igenväxande sandhed"
$ws.Range("AO11").ClearContents()
$ws.Range("AQ11").ClearContents()
$ws.Range("AR11").ClearContents()
$ws.Range("AX11").Value = "Nils Otto Nilsson, Mats Karlsson"
$ws.Range("A12").Value = 112154272
$ws.Range("B12").Value = 44328
$ws.Range("E12").Value = 102366
$ws.Range("F12").Value = "Ängsmetallvinge"
$ws.Range("G12").Value = "Adscita statices"
$ws.Range("H12").Value = "(Linnaeus, 1758)"
$ws.Range("I12").Value = "'3"
$ws.Range("M12").Value = "vilande"
$ws.Range("Q12").Value = 442480
$ws.Range("R12").Value = 6204371
$ws.Range("AO12").Value = "på blmr av åkervädd m fl"
$ws.Range("A13").Value = 112154283
$ws.Range("B13").Value = 44337
$ws.Range("E13").Value = 201164
$ws.Range("F13").Value = "Sexfläckig bastardsvärmare"
$ws.Range("G13").Value = "Zygaena filipendulae"
$ws.Range("I13").Value = "'1"
$ws.Range("M13").Value = "födosökande"
$ws.Range("Q13").Value = 442664
$ws.Range("R13").Value = 6204260
$ws.Range("AO13").Value = "på blmr av åkervädd"
$ws.Range("A14").Value = 112145591
$ws.Range("B14").Value = 42552
$ws.Range("E14").Value = 102923
$ws.Range("F14").Value = "Violettkantad guldvinge"
$ws.Range("G14").Value = "Lycaena hippothoe"
$ws.Range("H14").Value = "(Linnaeus, 1760)"
$ws.Range("K14").Value = "imago/adult"
$ws.Range("L14").Value = "hane"
$ws.Range("M14").Value = "födosökande"
$ws.Range("Q14").Value = 442543
$ws.Range("R14").Value = 6204460
$ws.Range("Y14").Value = "2013-06-12"
$ws.Range("AA14").Value = "2013-06-12"
$ws.Range("AC14").Value = "lufthåvning"
$ws.Range("AI14").Value = "i igenväxande hed"
$ws.Range("AO14").Value = "på tjärblomster"
$ws.Range("AQ14").Value = "Nils Otto Nilsson"
$ws.Range("AR14").Value = "NON 04542"
$ws.Range("AX14").Value = "Nils Otto Nilsson"
$ws.Range("A15").Value = 112154276
$ws.Range("B15").Value = 39455
$ws.Range("E15").Value = 102471
$ws.Range("F15").Value = "Åkerväddsantennmal"
$ws.Range("G15").Value = "Nemophora metallica"
$ws.Range("H15").Value = "(Poda, 1761)"
$ws.Range("I15").Value = "'1"
$ws.Range("Q15").Value = 442616
$ws.Range("R15").Value = 6204441
$ws.Range("AO15").Value = "på blmr av åkervädd"
$ws.Range("A16").Value = 112154275
$ws.Range("B16").Value = 44328
$ws.Range("E16").Value = 102366
$ws.Range("F16").Value = "Ängsmetallvinge"
$ws.Range("G16").Value = "Adscita statices"
$ws.Range("H16").Value = "(Linnaeus, 1758)"
$ws.Range("I16").Value = "'5"
$ws.Range("K16").ClearContents()
$ws.Range("L16").ClearContents()
$ws.Range("M16").Value = "vilande"
$ws.Range("Q16").Value = 442616
$ws.Range("R16").Value = 6204441
$ws.Range("Y16").Value = "2013-07-09"
$ws.Range("AA16").Value = "2013-07-09"
$ws.Range("AI16").Value = "på igenväxande sandhed"
$ws.Range("AO16").Value = "på blmr av åkervädd m fl"
$ws.Range("AX16").Value = "Nils Otto Nilsson, Mats Karlsson"
$ws.Range("A17").Value = 112154282
$ws.Range("B17").Value = 44328
$ws.Range("E17").Value = 102366
$ws.Range("F17").Value = "Ängsmetallvinge"
$ws.Range("G17").Value = "Adscita statices"
$ws.Range("I17").Value = "'4"
$ws.Range("M17").Value = "vilande"
